# Update "想去人数" (F column) counts across sheets to match the
# refreshed bilibili scrape output (gh-pages regeneration @ 456a3b4).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 161
$ws1.Range("F8").Value = 785
$ws1.Range("F14").Value = 5988
$ws1.Range("F16").Value = 2309
$ws1.Range("F18").Value = 163
$ws1.Range("F19").Value = 459
$ws1.Range("F20").Value = 9015
$ws1.Range("F22").Value = 2421
$ws1.Range("F24").Value = 2294
$ws1.Range("F25").Value = 2411
$ws1.Range("F26").Value = 1380
$ws1.Range("F27").Value = 232
$ws1.Range("F28").Value = 1941
$ws1.Range("F30").Value = 56
$ws1.Range("F40").Value = 93
$ws1.Range("F42").Value = 1510
$ws1.Range("F43").Value = 2461
$ws1.Range("F45").Value = 909
$ws1.Range("F46").Value = 289
$ws1.Range("F48").Value = 11
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F22").Value = 41
$ws2.Range("F23").Value = 41
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 687
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 687
$ws4.Range("F11").Value = 161
$ws4.Range("F12").Value = 785
$ws4.Range("F19").Value = 5988
$ws4.Range("F21").Value = 2309
$ws4.Range("F22").Value = 163
$ws4.Range("F23").Value = 459
$ws4.Range("F24").Value = 9015
$ws4.Range("F27").Value = 2422
$ws4.Range("F28").Value = 2294
$ws4.Range("F29").Value = 2411
$ws4.Range("F30").Value = 1380
$ws4.Range("F31").Value = 232
$ws4.Range("F32").Value = 1941
$ws4.Range("F34").Value = 56
$ws4.Range("F41").Value = 93
$ws4.Range("F43").Value = 1510
$ws4.Range("F44").Value = 2461
$ws4.Range("F45").Value = 909
$ws4.Range("F46").Value = 289
$ws4.Range("F51").Value = 41
